# Updated cryptos list (price/volume refresh), incl. a Polkadot/WrappedEther
# row swap at 13/14. D-column values are prefixed with a leading apostrophe
# so Excel stores the numeric-looking text (e.g. "93.60", "1.000") verbatim
# as text instead of silently converting it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.887.29"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "'1.888.67"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'0.7667"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "'242.78"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("D11").Value = "'0.08512"
$ws.Range("E11").Value = "  +4.63%  "
$ws.Range("D12").Value = "'0.7646"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.363"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.864.94"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'93.60"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'6.149"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "'29.849.37"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'13.78"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "'244.35"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "'0.000007805"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'2.143.12"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'7.984"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'0.1635"
$ws.Range("E25").Value = "  +3.95%  "
$ws.Range("D26").Value = "'9.425"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'162.47"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "'1.532"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").Value = "'4.498"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "'4.090"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").Value = "'0.7415"
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'2.696"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("D39").Value = "'0.01948"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'0.4468"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "'1.102.39"
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("D43").Value = "'72.96"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").Value = "'6.061"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "'0.8517"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'102.91"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "'1.867"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "'7.652"
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").Value = "'2.989"
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("D51").Value = "'2.043.08"
$ws.Range("E51").Value = "  -0.44%  "
